$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# Locate the body placeholder shape that holds the "Day N - ..." agenda text.
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $cand = $s.Shapes.Item($i)
    if ($cand.HasTextFrame -and $cand.TextFrame.HasText) {
        if ($cand.TextFrame.TextRange.Text -like "*Binding Track Activity*") {
            $shp = $cand
        }
    }
}

$tr = $shp.TextFrame.TextRange

# Find the paragraph that contains "Binding Track Activity" (the "Day 3 - ..." line).
$paraCount = $tr.Paragraphs().Count
$paraIndex = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $cand = $tr.Paragraphs($i, 1)
    if ($cand.Text -like "*Binding Track Activity*") {
        $paraIndex = $i
    }
}
$para = $tr.Paragraphs($paraIndex, 1)

# Find the run inside that paragraph holding the whole sentence.
$runIndex = -1
for ($i = 1; $i -le $para.Count; $i++) {
    $cand = $para.Runs($i, 1)
    if ($cand.Text -like "*Binding Track Activity*") {
        $runIndex = $i
    }
}

$run = $para.Runs($runIndex, 1)

# Split "Binding Track Activity (show / hide 'Loading' )" into three runs:
#   "Binding Track Activity (show / hide 'Loading"
#   "' "
#   "), Scan Operator"
$run.Text = "Binding Track Activity (show / hide " + [char]0x2018 + "Loading"
$null = $run.InsertAfter([char]0x2019 + " ")

$para = $tr.Paragraphs($paraIndex, 1)
$run2 = $para.Runs($runIndex + 1, 1)
$null = $run2.InsertAfter("), Scan Operator")
